# Update the "Pais" worksheet with the newest COVID-19 snapshot
# (commit: "Update countries & provincias Spain").
#
# The source data refresh reshuffled a few countries in the ranking
# (Barein/Israel, Islas Turcas y Caicos/Santa Sede and
# Islas Virgenes Britanicas/Papua Nueva Guinea swapped places) and
# brought several case counts up to date. The footer timestamp was
# also bumped to the new extraction time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer: "Datos actualizados a ..." timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 19 de Junio de 2020 a las 00:20"

# Row 4 - Estados Unidos (#8): refreshed totals
$ws.Range("B4").Value = 2258361
$ws.Range("C4").Value = 22634
$ws.Range("D4").Value = 925104
$ws.Range("E4").Value = 1212715
$ws.Range("G4").Value = 601
$ws.Range("H4").Value = 120542

# Row 5 - Brasil (#9): refreshed totals
$ws.Range("B5").Value = 978142
$ws.Range("C5").Value = 17833
$ws.Range("E5").Value = 426887
$ws.Range("G5").Value = 1083
$ws.Range("H5").Value = 47748

# Row 10 - Peru (#14): refreshed totals
$ws.Range("B10").Value = 244388
$ws.Range("C10").Value = 3480
$ws.Range("D10").Value = 131190
$ws.Range("E10").Value = 105737
$ws.Range("G10").Value = 204
$ws.Range("H10").Value = 7461

# Rows 50/51 - Barein and Israel swap ranking positions (#54/#55)
# and both receive refreshed totals.
$ws.Range("A50").Value = "Barein"
$ws.Range("B50").Value = 20430
$ws.Range("C50").Value = 469
$ws.Range("D50").Value = 14696
$ws.Range("E50").Value = 5679
$ws.Range("G50").Value = 6
$ws.Range("H50").Value = 55

$ws.Range("A51").Value = "Israel"
$ws.Range("B51").Value = 20036
$ws.Range("C51").Value = 253
$ws.Range("D51").Value = 15518
$ws.Range("E51").Value = 4215
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 303

# Row 54 - Japon (#58): refreshed totals
$ws.Range("B54").Value = 17668
$ws.Range("C54").Value = 40
$ws.Range("D54").Value = 15930
$ws.Range("E54").Value = 803
$ws.Range("G54").Value = 4
$ws.Range("H54").Value = 935

# Row 76 - Uzbekistan (#80): refreshed totals
$ws.Range("B76").Value = 5767
$ws.Range("C76").Value = 85
$ws.Range("E76").Value = 1582

# Row 90 - Bulgaria (#94): refreshed totals
$ws.Range("B90").Value = 3674
$ws.Range("C90").Value = 132
$ws.Range("D90").Value = 1941
$ws.Range("E90").Value = 1543
$ws.Range("G90").Value = 6
$ws.Range("H90").Value = 190

# Rows 208/209 - Islas Turcas y Caicos and Santa Sede swap positions
$ws.Range("A208").Value = "Islas Turcas y Caicos"
$ws.Range("D208").Value = 11
$ws.Range("H208").Value = 1

$ws.Range("A209").Value = "Santa Sede"
$ws.Range("D209").Value = 12
$ws.Range("H209").Value = 0

# Rows 213/214 - Islas Virgenes Britanicas and Papua Nueva Guinea swap positions
$ws.Range("A213").Value = "Islas Virgenes Britanicas"
$ws.Range("D213").Value = 7
$ws.Range("H213").Value = 1

$ws.Range("A214").Value = "Papua Nueva Guinea"
$ws.Range("D214").Value = 8
$ws.Range("H214").Value = 0
